# Update "想去人数" (F column) figures for a couple of events that were
# refreshed when the site was regenerated (gh-pages output at 456a3b4).
#
# Sheet "展览" (rId1 / sheet1.xml):
#   F2: 457  -> 458
#   F3: 5515 -> 5524
#   F4: 389  -> 388
#   F10: 13  -> 14
#
# Sheet "全部类型" (rId4 / sheet4.xml) mirrors the same rows (plus rows
# from the other sheets), so the same events need the same update there:
#   F2: 457  -> 458
#   F3: 5515 -> 5524
#   F4: 389  -> 388
#   F12: 13  -> 14

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 458
$wsExhibit.Range("F3").Value = 5524
$wsExhibit.Range("F4").Value = 388
$wsExhibit.Range("F10").Value = 14

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 458
$wsAll.Range("F3").Value = 5524
$wsAll.Range("F4").Value = 388
$wsAll.Range("F12").Value = 14
